$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values on specific rows to match repulled data / mean calc
$ws.Range("F3").Value = -11
$ws.Range("F4").Value = -10
$ws.Range("F6").Value = -4
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = -13
$ws.Range("F18").Value = -6
$ws.Range("F21").Value = 1
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = -2
